$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: new date value in BB1, copy formatting (style) from BA1 ---
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("BB1").Value = 45986

# --- Rows 2-70: BB column is an exact copy of the BA column values (no style) ---
$ws.Range("BA2:BA70").Copy()
$ws.Range("BB2:BB70").PasteSpecial(-4163)

# --- Rows 71-83: BB column gets new, updated forecast values ---
$ws.Range("BB71").Value = -0.1118837721692358
$ws.Range("BB72").Value = 0.3266766184601977
$ws.Range("BB73").Value = 0.325608361860148
$ws.Range("BB74").Value = 0.2086661300929905
$ws.Range("BB75").Value = 0.2086661300929905
$ws.Range("BB76").Value = 0.2086661300929905
$ws.Range("BB77").Value = 0.2086661300929905
$ws.Range("BB78").Value = 0.2086661300929905
$ws.Range("BB79").Value = 0.2086661300929905
$ws.Range("BB80").Value = 0.2086661300929905
$ws.Range("BB81").Value = 0.2086661300929905
$ws.Range("BB82").Value = 0.2086661300929905

# --- New row 83: date in column A (style copied from A82), and BB83 value ---
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.2086661300929905

$excel.CutCopyMode = 0
